$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 131273746
$ws.Range("B5").Value = 57884
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("M5").Value = "färska spår"
$ws.Range("P5").Value = "Sims fäbodar, Dlr"
$ws.Range("Q5").Value = 515255
$ws.Range("R5").Value = 6705041
$ws.Range("S5").Value = 50
$ws.Range("T5").Value = "Dalarna"
$ws.Range("U5").Value = "Borlänge"
$ws.Range("V5").Value = "Dalarna"
$ws.Range("W5").Value = "Stora Tuna"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2026-02-23"
$ws.Range("Y5").Style = "Normal"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2026-02-23"
$ws.Range("AA5").Style = "Normal"
$ws.Range("AC5").Value = "Ringhack på tall."
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = "Anna-Lena Thommson"
$ws.Range("AX5").Value = "Anna-Lena Thommson"

# Row 6
$ws.Range("A6").Value = 131273722
$ws.Range("B6").Value = 57884
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("M6").Value = "äldre spår"
$ws.Range("P6").Value = "Sims fäbodar, Dlr"
$ws.Range("Q6").Value = 515365
$ws.Range("R6").Value = 6705054
$ws.Range("S6").Value = 50
$ws.Range("T6").Value = "Dalarna"
$ws.Range("U6").Value = "Borlänge"
$ws.Range("V6").Value = "Dalarna"
$ws.Range("W6").Value = "Stora Tuna"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2026-02-23"
$ws.Range("Y6").Style = "Normal"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2026-02-23"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AC6").Value = "Ringhack på tall."
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = "Anna-Lena Thommson"
$ws.Range("AX6").Value = "Anna-Lena Thommson"
